# Auto-generated edit script: apply scheduled-runner price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1150.4117
$ws.Range("I15").Value = 1150.4117
$ws.Range("K15").Value = 3451.2351
$ws.Range("M15").Value = -3282.2351
$ws.Range("H132").Value = 2698.274
$ws.Range("I132").Value = 2725.2646
$ws.Range("K132").Value = 8175.793799999999
$ws.Range("M132").Value = -5645.793799999999
$ws.Range("H135").Value = 1060.6316
$ws.Range("I135").Value = 1060.6316
$ws.Range("K135").Value = 9545.6844
$ws.Range("M135").Value = -7010.6844
$ws.Range("H138").Value = 4188.9077
$ws.Range("J138").Value = 4753.65
$ws.Range("L138").Value = 14260.95
$ws.Range("N138").Value = -24540.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1610.9865
$ws.Range("I32").Value = 835.6957
$ws.Range("K32").Value = 835.6957
$ws.Range("M32").Value = -548.6957
$ws.Range("H74").Value = 1787.6154
$ws.Range("I74").Value = 1025.5385
$ws.Range("J74").Value = 2549.6924
$ws.Range("K74").Value = 1025.5385
$ws.Range("L74").Value = 2549.6924
$ws.Range("M74").Value = -151.5385000000001
$ws.Range("N74").Value = -4297.6924
$ws.Range("H77").Value = 1787.6154
$ws.Range("I77").Value = 1025.5385
$ws.Range("J77").Value = 2549.6924
$ws.Range("K77").Value = 5127.692500000001
$ws.Range("L77").Value = 12748.462
$ws.Range("M77").Value = -759.692500000001
$ws.Range("N77").Value = -21484.462
$ws.Range("H122").Value = 3086.5833
$ws.Range("I122").Value = 1966.6
$ws.Range("J122").Value = 3886.5715
$ws.Range("K122").Value = 5899.799999999999
$ws.Range("L122").Value = 11659.7145
$ws.Range("M122").Value = -3449.799999999999
$ws.Range("N122").Value = -16559.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 47995
$ws.Range("J6").Value = 47995
$ws.Range("L6").Value = 47995
$ws.Range("N6").Value = -48221
$ws.Range("H105").Value = 3533.0908
$ws.Range("I105").Value = 3533.0908
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3533.0908
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1786.0908
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 6955.9375
$ws.Range("I107").Value = 7648.0713
$ws.Range("J107").Value = 2111
$ws.Range("K107").Value = 7648.0713
$ws.Range("L107").Value = 2111
$ws.Range("M107").Value = -5728.0713
$ws.Range("N107").Value = -5951

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3011.0862
$ws.Range("I31").Value = 2224.0688
$ws.Range("J31").Value = 3798.1035
$ws.Range("K31").Value = 2224.0688
$ws.Range("L31").Value = 3798.1035
$ws.Range("M31").Value = -1929.0688
$ws.Range("N31").Value = -4388.1035
$ws.Range("H34").Value = 3011.0862
$ws.Range("I34").Value = 2224.0688
$ws.Range("J34").Value = 3798.1035
$ws.Range("K34").Value = 2224.0688
$ws.Range("L34").Value = 3798.1035
$ws.Range("M34").Value = -2022.0688
$ws.Range("N34").Value = -4202.1035
$ws.Range("H86").Value = 4951.25
$ws.Range("I86").Value = 4456.909
$ws.Range("J86").Value = 6038.8
$ws.Range("K86").Value = 4456.909
$ws.Range("L86").Value = 6038.8
$ws.Range("M86").Value = -3333.909
$ws.Range("N86").Value = -8284.799999999999
$ws.Range("H89").Value = 4951.25
$ws.Range("I89").Value = 4456.909
$ws.Range("J89").Value = 6038.8
$ws.Range("K89").Value = 22284.545
$ws.Range("L89").Value = 30194
$ws.Range("M89").Value = -16668.545
$ws.Range("N89").Value = -41426
$ws.Range("H107").Value = 336
$ws.Range("I107").Value = 327.88235
$ws.Range("J107").Value = 353.25
$ws.Range("K107").Value = 327.88235
$ws.Range("L107").Value = 353.25
$ws.Range("M107").Value = 1592.11765
$ws.Range("N107").Value = -4193.25
$ws.Range("H134").Value = 3778.7812
$ws.Range("I134").Value = 3135.3103
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 9405.930899999999
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -6870.930899999999
$ws.Range("N134").Value = -35067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 24741.676
$ws.Range("I141").Value = 6365.095
$ws.Range("K141").Value = 19095.285
$ws.Range("M141").Value = -13915.285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7887.1665
$ws.Range("I70").Value = 7638.6
$ws.Range("J70").Value = 8384.299999999999
$ws.Range("K70").Value = 7638.6
$ws.Range("L70").Value = 8384.299999999999
$ws.Range("M70").Value = -7368.6
$ws.Range("N70").Value = -8924.299999999999
$ws.Range("H73").Value = 7887.1665
$ws.Range("I73").Value = 7638.6
$ws.Range("J73").Value = 8384.299999999999
$ws.Range("K73").Value = 7638.6
$ws.Range("L73").Value = 8384.299999999999
$ws.Range("M73").Value = -6702.6
$ws.Range("N73").Value = -10256.3
$ws.Range("H80").Value = 2616.2727
$ws.Range("I80").Value = 2579.8333
$ws.Range("J80").Value = 2660
$ws.Range("K80").Value = 2579.8333
$ws.Range("L80").Value = 2660
$ws.Range("M80").Value = -1581.8333
$ws.Range("N80").Value = -4656
$ws.Range("H83").Value = 2616.2727
$ws.Range("I83").Value = 2579.8333
$ws.Range("J83").Value = 2660
$ws.Range("K83").Value = 12899.1665
$ws.Range("L83").Value = 13300
$ws.Range("M83").Value = -7907.166499999999
$ws.Range("N83").Value = -23284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1899.341
$ws.Range("I136").Value = 1455
$ws.Range("J136").Value = 2851.5
$ws.Range("K136").Value = 4365
$ws.Range("L136").Value = 8554.5
$ws.Range("M136").Value = -1815
$ws.Range("N136").Value = -13654.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 782.9375
$ws.Range("I107").Value = 602.5
$ws.Range("J107").Value = 1324.25
$ws.Range("K107").Value = 1807.5
$ws.Range("L107").Value = 3972.75
$ws.Range("M107").Value = 112.5
$ws.Range("N107").Value = -7812.75
$ws.Range("H132").Value = 3498.7673
$ws.Range("I132").Value = 3488.6487
$ws.Range("J132").Value = 3561.1667
$ws.Range("K132").Value = 10465.9461
$ws.Range("L132").Value = 10683.5001
$ws.Range("M132").Value = -7935.946100000001
$ws.Range("N132").Value = -15743.5001
$ws.Range("H136").Value = 2617.8928
$ws.Range("I136").Value = 2682.6365
$ws.Range("J136").Value = 2380.5
$ws.Range("K136").Value = 8047.9095
$ws.Range("L136").Value = 7141.5
$ws.Range("M136").Value = -5497.9095
$ws.Range("N136").Value = -12241.5

